$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing row 2 (Beta) values ---
$ws.Range("C2").Value = 19.12075701903682
$ws.Range("D2").Value = 0.006944391349456487
$ws.Range("E2").Value = 0.01982943797740053
$ws.Range("F2").Value = 6.94746946166665
$ws.Range("G2").Value = 5.138888063632246
$ws.Range("H2").Value = 9.401137236739489
$ws.Range("I2").Value = 0.00231147349769847
$ws.Range("J2").Value = 0.0008074615548642856
$ws.Range("K2").Value = 0.003346558989260584
$ws.Range("L2").Value = 0.009827919301656317
$ws.Range("M2").Value = 0.007836072250056482
$ws.Range("N2").Value = 0.01236402518280967

# --- Update existing row 3 (Gamma) values ---
$ws.Range("C3").Value = 0.04981522627320694
$ws.Range("D3").Value = 0.04815098319456564
$ws.Range("E3").Value = 0.0499839736740351
$ws.Range("F3").Value = 0.1355894329601918
$ws.Range("G3").Value = 0.002638118362494868
$ws.Range("H3").Value = 0.3231669579907591
$ws.Range("I3").Value = 0.125572212602547
$ws.Range("J3").Value = 0.002461095202072723
$ws.Range("K3").Value = 0.2987976006760024
$ws.Range("L3").Value = 0.1426215752832571
$ws.Range("M3").Value = 0.002754788459134767
$ws.Range("N3").Value = 0.3405286356767039

# --- Add new row 4 (Beta + Gamma) ---
$ws.Range("A4").Value = 2
# Copy the formatting (bold font + thin border + alignment) from A2 so the
# new "index" cell in column A matches the look of the other rows.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 19.17057224531003
$ws.Range("D4").Value = 0.05509537454402212
$ws.Range("E4").Value = 0.06981341165143562
$ws.Range("F4").Value = 7.083058894626841
$ws.Range("G4").Value = 5.14152618199474
$ws.Range("H4").Value = 9.724304194730248
$ws.Range("I4").Value = 0.1278836861002454
$ws.Range("J4").Value = 0.003268556756937009
$ws.Range("K4").Value = 0.3021441596652629
$ws.Range("L4").Value = 0.1524494945849134
$ws.Range("M4").Value = 0.01059086070919125
$ws.Range("N4").Value = 0.3528926608595136
